$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the TOTAL formula range (C4:C999 -> C4:C998)
$ws.Range("C2").Formula = "=SUM(C4:C998)"

# Copy the formatting of row 8 (A8:C8 -> style pattern s=4, s=15, s=14) onto the new row 10
$ws.Range("A8:C8").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Fill in the new TP6 group row
$ws.Range("A10").Value = "TP6"
$ws.Range("C11").Value = " "
$ws.Range("B10").Value = "Remodularização (os dois apresentaram)
Django monolítico para Restful (ideia desafiadora, parabéns)
MTV para Restful
AST, Astor
- Fizeram um exemplo simples mas que ilustra o processo
- Apontou as limitações. Importante!"
$ws.Range("C10").Value = 7.5
$ws.Rows.Item(10).RowHeight = 136

# Update the sheet view to match the new scroll/selection position
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("G8").Select()
